$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "https://www.ziprecruiter.com/k/l/AAK8EKdHGYONFjs4XxT4Dvhf8esdxhRqthjwYcb109juhf3CAg2HtrDnmMkb_FMIl3j409VM7Mv61aX5JYddyCFmaYDvGvLyCB43l-lSJDDrBKglhvhz0DZovG2e3iOpnEVfTa_LibHS5hLv_TxXOaa-pU6V3WXP6Mf0An7LEzliuiGHXsM61g"
$ws.Range("F3").Value = "https://www.ziprecruiter.com/k/l/AAK-Jb_bcveOnTcdYOCWuYeewTFCTtHJtTPEKD4ewlIeWTddR9qGe3G0LGGdnB4jTu0kY6ySX2XrRW-PO1YCNqJ7yeDfW1VIwM8KVNIxaPRbjqjkI1qt7calYtYzhqbpIbdfd4DU83nwkvv3Zj-MqKNFWoUwbf9NAlS9Pl4hfI4TXFvVjHaCYg"
$ws.Range("F4").Value = "https://www.ziprecruiter.com/k/l/AAJ0z_Ri8IEIA96X8ACCeZelvcn2mbH5h24TWfMcQJVgJFdwfOq6vkX0Pi-C1APb2fP55T7rY4I0uunWrdM9k8BhKYe09jY8RulE6eRQJ9sBHO76XFtarg_495_g37UnR58bJ3fVaChk8uQ0XvDTdEYVW1QcdVqnlCsSdGo1VJXAcBYJVcfVyA"
$ws.Range("F5").Value = "https://www.ziprecruiter.com/k/l/AAIaXYqaVYgwohWgTlQsZWkJclWLNOh8PYIRAFr0xUORBtFl0iCdqc-GpTta9hpqv2B6Zabmf9Jn-H80jKVUBsMNm9x25fLHopI0xh-5lK977mu9laWss9ckUkAVXatO4YGpgsQK0Thmdwd1rLNFQFeogxYjYvRcadc11fAw6ioVNzHz9A_44g"
$ws.Range("F6").Value = "https://www.ziprecruiter.com/ek/l/AAIZ15FtHnM1nCAldE0gUjmfWXD2r-QaYDT8o2bp2BMlv2VDA0IWKdWDltfMOPfeZr2VkjDknzUKG34IbrEP_cHXllTsT6hrmAKxODOeifjTdLFMIU4MRoEXSECplm4vTfvAQlehjFikKDWiQxxScKsVR7O2GDbiYaseMRHZdPPpKSyopQOgig"
$ws.Range("F7").Value = "https://www.ziprecruiter.com/k/l/AAKgDsqtHLCbsfEcjRdJQHPKA8-Yv904kk_paKPKWfwLZyzMRG9shPKvErBtP5YeBVIXvOGRpsBeLUhmGbNcXdq21vzOM6ALrsejlCh55h2yRhql0d4pUIz8NJ7XDWjgXgQZA-3hZSnXAzCHlW3USBpCJ6M817sb7jFi0cSaUHOPdU29cAW8BA"
$ws.Range("F8").Value = "https://www.ziprecruiter.com/k/l/AAKoPjQiU7W7jRnszW0KEG06BGfySCHLOaYl_PF2eI3Y7FVm-YSSk0sRZR8CfBdcYffb4we8ytqpKY3W01YtG8bSS3iTzgAL7pr-fZQvfg_7FXHDgV8v5OvftuQ397csEdWnuNkj8Cs1K0NrM11gkfcercc-ZdL1XbP4fBAG-jOyLhKK9sOI7g"
$ws.Range("F9").Value = "https://www.ziprecruiter.com/k/l/AAK2M4cxQkrEl38mxj6YOe42G8lUqJcotFvFVli83baTm30QAfAvm_YIummYxbR0WaUkaY6B4SQkl9KTr9TQ8DkZHvBvMFwUIbSJEqKxZW1OW7WDvgvO6zZyTRAdNmW0RqwkdP64yQO7Mez9_8DToRxRzgSpxgeFd3wRJ3jqo0RTSX7-gJzBqg"
$ws.Range("F10").Value = "https://www.ziprecruiter.com/k/l/AAL48Rh1EHGsBJq_5gzCQmmqb6GxQxpNG04GUxfZdD9AtKupDLzFNpe4LTInK51vJcvRsCof8w-aHCpgiMxYco1ErSvUVs5_-lPmworZtVrm-BJau4iUEkzxIb1JIrgA883XGCbw1jn59LxR_p-lKJYSUbp4Rhv9q6D1bQjJRQ0bLKI3k4llyA"
$ws.Range("F11").Value = "https://www.ziprecruiter.com/k/l/AAIBquUeYBUSZL0QusytWw5q7xQt6gbpfpQaDqyaDEsP39kHYOL_6Qtzbg-13hFesDe77A65xGJmYA2I9y5_QeY9xX44FJquOUvP9J7VIA54DVy6niQBtQ0YkSsD6jpkoKCvZ4_CjbnzvjR1hzzJe5vBRK60AmF0Zt8QVGCQbK3CRvYp3V_eXmA"
